# Apply updated dSF (column F) values, re-pulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    14 = -1
    29 = 0
    30 = 1
    31 = 7
    33 = 0
    35 = -8
    42 = 0
    44 = -3
    45 = -1
    54 = 4
    55 = -6
    58 = -8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
